{"js": "// Replace the date line and all 25 two-digit multiplication answers with\n// the new values from the latest daily-practice run.\nconst replacements = [\n  [\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"],\n  [\"38\u00d792=3496\", \"14\u00d779=1106\"],\n  [\"67\u00d727=1809\", \"95\u00d747=4465\"],\n  [\"68\u00d799=6732\", \"33\u00d756=1848\"],\n  [\"29\u00d737=1073\", \"51\u00d718=918\"],\n  [\"60\u00d782=4920\", \"75\u00d724=1800\"],\n  [\"16\u00d756=896\", \"25\u00d799=2475\"],\n  [\"82\u00d722=1804\", \"55\u00d757=3135\"],\n  [\"65\u00d798=6370\", \"26\u00d761=1586\"],\n  [\"44\u00d783=3652\", \"39\u00d797=3783\"],\n  [\"75\u00d725=1875\", \"11\u00d720=220\"],\n  [\"58\u00d784=4872\", \"85\u00d778=6630\"],\n  [\"78\u00d742=3276\", \"54\u00d765=3510\"],\n  [\"55\u00d723=1265\", \"42\u00d738=1596\"],\n  [\"76\u00d718=1368\", \"42\u00d784=3528\"],\n  [\"87\u00d787=7569\", \"38\u00d723=874\"],\n  [\"23\u00d722=506\", \"59\u00d783=4897\"],\n  [\"36\u00d757=2052\", \"33\u00d769=2277\"],\n  [\"47\u00d788=4136\", \"47\u00d787=4089\"],\n  [\"55\u00d772=3960\", \"53\u00d740=2120\"],\n  [\"56\u00d789=4984\", \"15\u00d772=1080\"],\n  [\"21\u00d747=987\", \"23\u00d773=1679\"],\n  [\"53\u00d747=2491\", \"25\u00d728=700\"],\n  [\"45\u00d718=810\", \"32\u00d761=1952\"],\n  [\"16\u00d733=528\", \"43\u00d736=1548\"],\n  [\"72\u00d729=2088\", \"16\u00d717=272\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 25 two-digit multiplication answers with\n# the new values from the latest daily-practice run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old='2025-09-22 Monday'; new='2025-09-23 Tuesday'},\n    @{old='38\u00d792=3496'; new='14\u00d779=1106'},\n    @{old='67\u00d727=1809'; new='95\u00d747=4465'},\n    @{old='68\u00d799=6732'; new='33\u00d756=1848'},\n    @{old='29\u00d737=1073'; new='51\u00d718=918'},\n    @{old='60\u00d782=4920'; new='75\u00d724=1800'},\n    @{old='16\u00d756=896'; new='25\u00d799=2475'},\n    @{old='82\u00d722=1804'; new='55\u00d757=3135'},\n    @{old='65\u00d798=6370'; new='26\u00d761=1586'},\n    @{old='44\u00d783=3652'; new='39\u00d797=3783'},\n    @{old='75\u00d725=1875'; new='11\u00d720=220'},\n    @{old='58\u00d784=4872'; new='85\u00d778=6630'},\n    @{old='78\u00d742=3276'; new='54\u00d765=3510'},\n    @{old='55\u00d723=1265'; new='42\u00d738=1596'},\n    @{old='76\u00d718=1368'; new='42\u00d784=3528'},\n    @{old='87\u00d787=7569'; new='38\u00d723=874'},\n    @{old='23\u00d722=506'; new='59\u00d783=4897'},\n    @{old='36\u00d757=2052'; new='33\u00d769=2277'},\n    @{old='47\u00d788=4136'; new='47\u00d787=4089'},\n    @{old='55\u00d772=3960'; new='53\u00d740=2120'},\n    @{old='56\u00d789=4984'; new='15\u00d772=1080'},\n    @{old='21\u00d747=987'; new='23\u00d773=1679'},\n    @{old='53\u00d747=2491'; new='25\u00d728=700'},\n    @{old='45\u00d718=810'; new='32\u00d761=1952'},\n    @{old='16\u00d733=528'; new='43\u00d736=1548'},\n    @{old='72\u00d729=2088'; new='16\u00d717=272'}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.new\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
